$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 526, shifting rows 526:634 down to 527:635
$ws.Rows("526:526").Insert()

# Populate the new row 526 with values. Most columns mirror the row that is
# now at 527 (the former row 526), except for D, M, P and S which take new
# values as specified by the edit.
$ws.Range("A526").Value = 9
$ws.Range("B526").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C526").Value = "Metropolitana"
$ws.Range("D526").Value = 45015
$ws.Range("D526").NumberFormat = $ws.Range("D527").NumberFormat
$ws.Range("E526").Value = 13
$ws.Range("F526").Value = "Fruta"
$ws.Range("G526").Value = 100108
$ws.Range("H526").Value = "Tropicales y subtropicales"
$ws.Range("I526").Value = 100108002
$ws.Range("J526").Value = "Mango"
$ws.Range("K526").Value = "Sin especificar"
$ws.Range("L526").Value = "Primera"
$ws.Range("M526").Value = 580
$ws.Range("N526").Value = 6000
$ws.Range("O526").Value = 6500
$ws.Range("P526").Value = 6259
$ws.Range("Q526").Value = "`$/bandeja 4 kilos"
$ws.Range("R526").Value = "Perú"
$ws.Range("S526").Value = 1565
$ws.Range("T526").Value = 4
